$wb = $excel.ActiveWorkbook

# ----- Sheet2 (Ownership): update column widths -----
$wsOwnership = $wb.Worksheets.Item("Ownership")
$wsOwnership.Columns.Item(1).ColumnWidth = 58 - (5/6)
$wsOwnership.Columns.Item(2).ColumnWidth = 42.33203125 - (5/6)
$wsOwnership.Columns.Item(3).ColumnWidth = 42 - (5/6)

# ----- New content strings for row 4 -----
$task = @'
Microsoft was using Tenant Data store for Office 365 Admin Center Webapi's.
My Team Lead picked up that work item and it was time sensitive task and
Microsoft was retiring Tenant Data store in next two months. After a month my TeamLead
had to leave the project as he got another Full time project. I have taken this task 
to get it complete in a month based upon intial analysis made by my Team Lead.
Based upon intial analysis I came know that I have to migrate from Tenant Data store to
Azure Cosmos DB. Make API code changes to use Azure Cosmos DB and Remove Tenant Data store.
'@

$taskSteps = @'
I created some of the below steps:
step1: Create a POC to get connect to azure cosmos db collection and perform
       Update, Read/Write and Delete operations.
step2: Migrate existing data from TDS to AZCDB
step3: Make API changes to use Azure Cosmos DB instead of Tenant Data store.
step4: Apply POC implementation on API.
step5: Update unit tests and integration tests.
Step 6: Demo and Get approval.
step 7: Deploy it Production and make it available for 10% users using A/B test methods.

'@

$action = @'
Action:
1. At that time I was new to perform CRUD operations on Azure Cosmos DB.
   I created a POC to perform READ/Write opearation on Azure Cosmos DB test collection.
   I used Azure SDK to get connect to Azure Cosmos DB collection and Used Azure Libraries to
   Insert, Update, Delete to Azure Cosmos DB collection.
   I demoed it to the Team and Team accepted my POC. 
2. Migrate the data from Tenant Data store(Key Value pair data) to Azure Cosmos DB.
3. Web api's using Tenant data store to  read/write/update/delete the data. Made
   api code changes to read/write from Tenant data store to Azure Cosmos DB.
4. Based upon POC, I applied that implementation to actual API's.
5. Deployed the code to Test environments to give a demo Tenant Data Store vs Azure Cosmos DB.
5. Gave a demo of my code changes and team was happy with my code changes.

'@

$result = @'
Result:
The Api ran successfully as it ran on Tenant data store, after migrating to Azure Cosmos DB.
I deployed the old code to Test01 server and new code to Test02 environments. Compared the results
on both servers and demo the same to the team. Team was happy with the implementation. Team approved my changes to get it deployed
to production to get it released to 10% of the website users and monitor for few days. After monitoring,
No issues were reported and then made it available to 100% users.  
'@

# ----- Add new row 4 with Situation/Task/Action/Result content -----
$wsOwnership.Range("A4").Value = $task
$wsOwnership.Range("B4").Value = $taskSteps
$wsOwnership.Range("C4").Value = $action
$wsOwnership.Range("D4").Value = $result
$wsOwnership.Range("A4:D4").WrapText = $true
$wsOwnership.Rows.Item(4).RowHeight = 217.8

# ----- Activate Ownership sheet and select A4 (becomes the active sheet/tab) -----
$wsOwnership.Activate()
$wsOwnership.Range("A4").Select()
